$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FUELS")

# Update PEN (C2) for Natural Gas (NG) row - now computed from ecoinvent 3.4 components
$ws.Range("C2").Formula = "=1.1767+0.0019487+0.0000015726"

# Update CO2 (D2) for Natural Gas (NG) row
$ws.Range("D2").Value = 0.06682

# Update reference/source text (F2) to point to the new ecoinvent source description
$ws.Range("F2").Value = "ecoinvent 3.4 - market for natural gas, burned in gas motor, for storage_GLO_2017_Allocation, cut-off"

# Reflect the new active selection on the sheet
$ws.Range("C2").Select()
